# Add a new "Зарегистрирован" (Registered) boolean column and a new user row,
# matching the "added method recording in excel" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Пользователи")

# New header for column I
$ws.Range("I1").Value = "Зарегистрирован"
$ws.Columns.Item(9).ColumnWidth = 16.7

# Existing row 2 gets a value in the new column
$ws.Range("I2").Value = $true

# New user data in row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Марина"
$ws.Range("C3").Value = "Шарапова"
$ws.Range("D3").Value = "mail@mail.ru"
$ws.Range("F3").Value = "fsdfsdfs"
$ws.Range("E3").Value = "u3'*OlGq"
$ws.Range("G3").NumberFormat = "m/d/yy"
$ws.Range("G3").Value = 7337
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = $false

# New hyperlink for the e-mail address in the new row
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:mail@mail.ru")
$ws.Range("D3").Style = "Гиперссылка"

# Move the active selection, as recorded in the workbook
$ws.Range("H9").Select()
